# AP160_TestData_ManageAcctPeriods_21C.xlsx - "Add files via upload" edit
#
# The author re-uploaded this test-data workbook after scrubbing the
# environment-specific login fields on the Input_Value sheet: the URL,
# UserName and Password cells (L2:N2) are blanked out (their shared
# strings - the Oracle Cloud URL, the implementation-user login and the
# password - are removed from the workbook entirely), while keeping each
# cell's existing number/style formatting intact. All other data
# (wait-time constants, ledger/period search values, Output_Value and
# Compare sheets) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Make sure we're working on the right sheet (it's already the active one
# in the source file, but activating keeps tabSelected/the active-sheet
# state correct regardless of call order).
$ws.Activate()

# Clear the stored login values but keep cell formatting (style) as-is.
$loginRange = $ws.Range("L2:N2")
$loginRange.ClearContents()

# Reflect the real edit's on-screen selection (the three cells that were
# just cleared).
$loginRange.Select()
